# "Generate Report for Archive"
#
# The localization status text "Ready for handoff" is updated to
# "In Translation" everywhere it appears (Overview!E2:F4, zh-cn!C2:C4,
# de-de!C2:C4), and the status columns are re-sized to fit the new
# (shorter) text, shrinking from 17.2159881591797 down to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: column C ("Status") holds the status text ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: column C ("Status") holds the status text ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
